$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = "MAP_OPT.LZ"
$ws.Range("B14").Value = 436
$ws.Range("C14").Value = "mission options frame"

# Row 15
$ws.Range("A15").Value = "RAM2VERT.DTA"
$ws.Range("B15").Value = 10
$ws.Range("C15").Value = "frame vert part"

# Row 16
$ws.Range("A16").Value = "RAM2HORZ.DTA"
$ws.Range("B16").Value = 76
$ws.Range("C16").Value = "frame horz part"

# Row 17
$ws.Range("A17").Value = "RAM2ROH.DTA"
$ws.Range("B17").Value = 10
$ws.Range("C17").Value = "frame corner"

# Update selection to match diff
$ws.Range("K15").Select()
